$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: now the ORG_Website_Cache_Scenario row (moved up from row 4), with
#     the record-count bumped from 1 to 150 and Execution switched to Yes ---
$ws.Range("A2").Value = "0000_ORG_Website_Cache_Scenario"
$ws.Range("B2").Value = "Verify the ORG_Website_Cache_ Info"
$ws.Range("C2").Value = "Verify_ORG_WEBSITE_CACHE"
$ws.Range("D2").Value = "150"
$ws.Range("E2").Value = "Yes"

# --- Row 3: now the ORG_Phone_Scenario row (moved up from row 2) ---
$ws.Range("A3").Value = "0000_ORG_Phone_Scenario"
$ws.Range("B3").Value = "Verify the ORG_Phone Info"
$ws.Range("C3").Value = "Verify_All_Buckets_ORG_PHONE"
$ws.Range("D3").Value = "5"
$ws.Range("E3").Value = "No"

# --- Row 4: now the ORG_PV_Phone_Scenario row (moved down from row 3), with
#     Execution switched to No ---
$ws.Range("A4").Value = "0000_ORG_PV_Phone_Scenario"
$ws.Range("B4").Value = "Verify the ORG_PV_Phone Info"
$ws.Range("C4").Value = "Verify_All_Buckets_ORG_PV_PHONE"
$ws.Range("D4").Value = "100"
$ws.Range("E4").Value = "No"

# Widen column A to fit the longer scenario names
$ws.Columns("A").ColumnWidth = 49.67

# Move the active selection to K27, matching the saved view state
$ws.Range("K27").Select()
